$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CourtHearings")
$ws.Activate()

# Highlight row 9 (A9:AG9) with green fill, matching the rest of the row's style
$ws.Range("A9:AG9").Interior.Color = 5296274

# Add a new row 12, matching the yellow-highlight formatting used by the
# surrounding data rows (same fill as row 11: A11:AG11)
$ws.Range("A12:AG12").Interior.Color = 65535

$ws.Range("A12").Value = "testT4149"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 6
$ws.Range("D12").Value = "n/a"
$ws.Range("E12").Value = "Yes"
$ws.Range("F12").Value = "n/a"
$ws.Range("G12").Value = "n/a"
$ws.Range("H12").Value = "future"
$ws.Range("I12").Value = "n/a"
$ws.Range("J12").Value = "Requested"
$ws.Range("K12").Value = "n/a"
$ws.Range("L12").Value = "Click"
$ws.Range("M12").Value = "n/a"
$ws.Range("N12").Value = "n/a"
$ws.Range("O12").Value = "n/a"
$ws.Range("P12").Value = "n/a"
$ws.Range("Q12").Value = "n/a"
$ws.Range("R12").Value = "n/a"
$ws.Range("S12").Value = "n/a"
$ws.Range("T12").Value = "n/a"
$ws.Range("U12").Value = "n/a"
$ws.Range("V12").Value = "n/a"
$ws.Range("W12").Value = "n/a"
$ws.Range("X12").Value = "n/a"
$ws.Range("Y12").Value = "n/a"
$ws.Range("Z12").Value = "n/a"
$ws.Range("AA12").Value = "n/a"
$ws.Range("AB12").Value = "n/a"
$ws.Range("AC12").Value = "n/a"
$ws.Range("AD12").Value = "n/a"
$ws.Range("AE12").Value = "n/a"
$ws.Range("AF12").Value = "n/a"
$ws.Range("AG12").Value = "n/a"

$ws.Range("A12").Select()
